# Daily attendance processing - 2025-10-30 20:23:48
#
# Column G ("Recorded By") holds a comma-separated list of whoever
# recorded/touched a session's attendance. This pass rotates each
# multi-author list left by one position (the earliest entry moves to
# the back of the list), except for the canonical
# "backup@backdoor.com, System" pairing, which is left as-is.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

$skipValue = "backup@backdoor.com, System"

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value2

    if ($val -eq $null) { continue }
    if ($val -eq "") { continue }
    if ($val -eq $skipValue) { continue }

    $parts = $val -split ", "
    if ($parts.Length -le 1) { continue }

    $rotated = ($parts[1..($parts.Length - 1)] + $parts[0]) -join ", "
    $cell.Value = $rotated
}
